# Initial notes on key points discussed
#
# 1. Slide 6 ("Key points discussed") gets a new Content Placeholder
#    (idx=1) with three bullet points - cloned from the matching
#    placeholder on slide 8 ("Challenges") so it inherits the same
#    nvSpPr/cNvSpPr/nvPr/spPr/bodyPr/lstStyle skeleton, then the text is
#    filled in (third bullet has a bold "The challenge" run).
# 2. Slide 8's first bullet "Time" becomes "Time flies".

$p = $ppt.ActivePresentation

# --- Slide 6: add the new Content Placeholder with notes -------------
$s6 = $p.Slides.Item(6)
$s8 = $p.Slides.Item(8)

$srcPlaceholder = $s8.Shapes.Item(2)
[void]$srcPlaceholder.Copy()
[void]$s6.Shapes.Paste()
$newShape = $s6.Shapes.Item($s6.Shapes.Count)

$bullet1 = "The capabilities of SimSurvey and sdmTMB"
$bullet2 = "Tractable topics to explore"
$bullet3 = "Sub-group specific discussions of the focus topics listed in The challenge slide"

$newShape.TextFrame.TextRange.Text = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3

# Bold just the "The challenge" phrase inside the third bullet.
$boldPhrase = "The challenge"
$p3Start = $bullet1.Length + 1 + $bullet2.Length + 1 + $bullet3.IndexOf($boldPhrase) + 1
$newShape.TextFrame.TextRange.Characters($p3Start, $boldPhrase.Length).Font.Bold = $true

# --- Slide 8: "Time" -> "Time flies" ----------------------------------
$s8Body = $s8.Shapes.Item(2).TextFrame.TextRange
$s8Body.Characters(1, 4).Text = "Time flies"
